# Commit: "Updated Argent prices in Excel"
# Appends a new row (row 15) to every price sheet with date 2025-03-16 and
# that day's price, repeating the last known (2025-03-15) price value -
# matching the source diff exactly.

$wb = $excel.ActiveWorkbook

# Sheet name -> new price value for 2025-03-16 (row 15, column B)
$updates = [ordered]@{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.298"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,455"
    "Silver Busbar front-side"   = "8,167"
    "Silver finger front-side"   = "8,217"
    "USD_CNY"                    = "7.2637"
}

foreach ($name in $updates.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $price = $updates[$name]

    # Force text entry (the source column stores these as text, not dates /
    # numbers) by pre-formatting as Text before assigning the value.
    $ws.Range("A15").NumberFormat = "@"
    $ws.Range("A15").Value = "2025-03-16"

    $ws.Range("B15").NumberFormat = "@"
    $ws.Range("B15").Value = $price

    # Row 14 (2025-03-15) carries the plain "General"/default style used
    # throughout the sheet. Copy its formatting onto the new row so we don't
    # leave a stray Text-number-format style behind (matches the original
    # file, where every cell uses the default style).
    $ws.Range("A14:B14").Copy()
    $ws.Range("A15:B15").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
